$wb = $excel.ActiveWorkbook

# --- Pest_list: update pest name, zoom, and selection ---
$wsPest = $wb.Worksheets.Item("Pest_list")
$wsPest.Range("A2").Value = "Conotrachelus nenuphar"

# --- Other settings: update a few dropdown-backed values ---
$wsOther = $wb.Worksheets.Item("Other settings")
$wsOther.Range("B2").Value = "no"
$wsOther.Range("B3").Value = "USA"
$wsOther.Range("B6").Value = "no"

# --- tech: update numeric values for the USA row ---
$wsTech = $wb.Worksheets.Item("tech")
$wsTech.Range("C9").Value = -50
$wsTech.Range("E9").Value = 65

# --- View state: move the active tab from "Other settings" to "tech",
#     set Pest_list's zoom/selection, and update each sheet's selection ---
$wsPest.Activate()
$excel.ActiveWindow.Zoom = 85
$wsPest.Range("D9").Select()

$wsOther.Activate()
$wsOther.Range("B3").Select()

$wsTech.Activate()
$wsTech.Range("E9").Select()
